# Actualización automática desde WSL
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the existing row 18 timestamp (column A) - tiny fractional fix
$ws.Range("A18").Value = 45877.75021111111

# Append new row 19 with the new weather-station reading
$ws.Range("A19").Value = 45877.79190259425
$ws.Range("B19").Value = 2025
$ws.Range("C19").Value = 32
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 86.3
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 6.61
$ws.Range("H19").Value = "ENE"
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = "19:00:20"

# Column A uses the same date/time style as the rows above it
$ws.Range("A19").NumberFormat = $ws.Range("A18").NumberFormat
